# Add "2022-Q3" sheet (fund holdings) right after "总计", and update the
# "总计" (Total) summary sheet with a new leading row for 2022-Q3,
# shifting the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2    = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet by cloning "2022-Q2" (keeps styles) ---
$q2.Copy($null, $total)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# "2022-Q2" source only has 3 data rows (rows 2-4); the new sheet needs 5
# data rows (rows 2-6). Extend formatting from row 4 down into rows 5-6.
$q3.Range("A4:H4").Copy()
$q3.Range("A5:H6").PasteSpecial(-4122)

# Helper-less direct cell writes (row, col): 1=A,2=B,3=C,4=D,5=E,6=F,7=G,8=H

# Row 2: 515210 / 国泰中证钢铁ETF
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).NumberFormat = "@"
$q3.Cells.Item(2,2).Value = "515210"
$q3.Cells.Item(2,3).NumberFormat = "@"
$q3.Cells.Item(2,3).Value = "国泰中证钢铁ETF"
$q3.Cells.Item(2,4).NumberFormat = "@"
$q3.Cells.Item(2,4).Value = "14.23"
$q3.Cells.Item(2,5).NumberFormat = "@"
$q3.Cells.Item(2,5).Value = "97.88"
$q3.Cells.Item(2,6).NumberFormat = "@"
$q3.Cells.Item(2,6).Value = "2.81"
$q3.Cells.Item(2,7).NumberFormat = "@"
$q3.Cells.Item(2,7).Value = "0.3999"
$q3.Cells.Item(2,8).Value = 9

# Row 3: 163110 / 申万菱信量化小盘股票（LOF）A
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).NumberFormat = "@"
$q3.Cells.Item(3,2).Value = "163110"
$q3.Cells.Item(3,3).NumberFormat = "@"
$q3.Cells.Item(3,3).Value = "申万菱信量化小盘股票（LOF）A"
$q3.Cells.Item(3,4).NumberFormat = "@"
$q3.Cells.Item(3,4).Value = "5.04"
$q3.Cells.Item(3,5).NumberFormat = "@"
$q3.Cells.Item(3,5).Value = "93.06"
$q3.Cells.Item(3,6).NumberFormat = "@"
$q3.Cells.Item(3,6).Value = "0.59"
$q3.Cells.Item(3,7).NumberFormat = "@"
$q3.Cells.Item(3,7).Value = "0.0297"
$q3.Cells.Item(3,8).Value = 9

# Row 4: 012977 / 瑞达鑫红量化6个月持有混合A
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).NumberFormat = "@"
$q3.Cells.Item(4,2).Value = "012977"
$q3.Cells.Item(4,3).NumberFormat = "@"
$q3.Cells.Item(4,3).Value = "瑞达鑫红量化6个月持有混合A"
$q3.Cells.Item(4,4).NumberFormat = "@"
$q3.Cells.Item(4,4).Value = "0.43"
$q3.Cells.Item(4,5).NumberFormat = "@"
$q3.Cells.Item(4,5).Value = "94.69"
$q3.Cells.Item(4,6).NumberFormat = "@"
$q3.Cells.Item(4,6).Value = "4.93"
$q3.Cells.Item(4,7).NumberFormat = "@"
$q3.Cells.Item(4,7).Value = "0.0212"
$q3.Cells.Item(4,8).Value = 2

# Row 5: 012978 / 瑞达鑫红量化6个月持有混合C
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).NumberFormat = "@"
$q3.Cells.Item(5,2).Value = "012978"
$q3.Cells.Item(5,3).NumberFormat = "@"
$q3.Cells.Item(5,3).Value = "瑞达鑫红量化6个月持有混合C"
$q3.Cells.Item(5,4).NumberFormat = "@"
$q3.Cells.Item(5,4).Value = "0.11"
$q3.Cells.Item(5,5).NumberFormat = "@"
$q3.Cells.Item(5,5).Value = "94.69"
$q3.Cells.Item(5,6).NumberFormat = "@"
$q3.Cells.Item(5,6).Value = "4.93"
$q3.Cells.Item(5,7).NumberFormat = "@"
$q3.Cells.Item(5,7).Value = "0.0054"
$q3.Cells.Item(5,8).Value = 2

# Row 6: 013918 / 申万菱信量化小盘股票（LOF）C
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).NumberFormat = "@"
$q3.Cells.Item(6,2).Value = "013918"
$q3.Cells.Item(6,3).NumberFormat = "@"
$q3.Cells.Item(6,3).Value = "申万菱信量化小盘股票（LOF）C"
$q3.Cells.Item(6,4).NumberFormat = "@"
$q3.Cells.Item(6,4).Value = "0.00"
$q3.Cells.Item(6,5).NumberFormat = "@"
$q3.Cells.Item(6,5).Value = "93.06"
$q3.Cells.Item(6,6).NumberFormat = "@"
$q3.Cells.Item(6,6).Value = "0.59"
$q3.Cells.Item(6,7).Value = 0
$q3.Cells.Item(6,8).Value = 9

# --- 2. Update the "总计" (Total) summary sheet ---
# Extend formatting (column A's bordered style) down from row 5 into row 6.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

# New leading row 2: 2022-Q3 — existing rows shift down by one.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 5
$total.Cells.Item(2,4).Value = 0.46

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q2"
$total.Cells.Item(3,3).Value = 3
$total.Cells.Item(3,4).Value = 0.57

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q1"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q4"
$total.Cells.Item(5,3).Value = 2
$total.Cells.Item(5,4).Value = 0.12

$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(6,2).Value = "2021-Q3"
$total.Cells.Item(6,3).Value = 9
$total.Cells.Item(6,4).Value = 8.2
